$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column cells we touch to stay text, matching the
# original inlineStr cells (Excel would otherwise auto-coerce values
# like "1.008" or "19.64" into numbers).
$priceCells = @("D2","D3","D4","D5","D7","D8","D9","D10","D11","D12","D13","D14","D15","D16","D17","D18","D19","D20","D21","D22","D23","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D39","D40","D41","D42","D43","D44","D47","D48","D49","D50","D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.930.29"
$ws.Range("E2").Value = "  -1.53%  "
$ws.Range("D3").Value = "1.831.62"
$ws.Range("E3").Value = "  -1.41%  "
$ws.Range("D4").Value = "1.008"
$ws.Range("E4").Value = "  +0.77%  "
$ws.Range("D5").Value = "310.93"
$ws.Range("E5").Value = "  -1.02%  "
$ws.Range("D7").Value = "0.4580"
$ws.Range("E7").Value = "  -0.77%  "
$ws.Range("D8").Value = "0.3692"
$ws.Range("E8").Value = "  -0.26%  "
$ws.Range("D9").Value = "0.07183"
$ws.Range("E9").Value = "  -1.95%  "
$ws.Range("D10").Value = "0.8784"
$ws.Range("E10").Value = "  -0.50%  "
$ws.Range("D11").Value = "0.07825"
$ws.Range("E11").Value = "  -0.36%  "
$ws.Range("D12").Value = "19.64"
$ws.Range("E12").Value = "  -0.98%  "
$ws.Range("D13").Value = "1.826.10"
$ws.Range("E13").Value = "  -2.26%  "
$ws.Range("D14").Value = "5.339"
$ws.Range("E14").Value = "  -0.78%  "
$ws.Range("D15").Value = "6.395"
$ws.Range("E15").Value = "  -2.45%  "
$ws.Range("D16").Value = "87.20"
$ws.Range("E16").Value = "  -5.04%  "
$ws.Range("D17").Value = "1.009"
$ws.Range("E17").Value = "  +0.73%  "
$ws.Range("D18").Value = "0.000008719"
$ws.Range("E18").Value = "  -1.51%  "
$ws.Range("D19").Value = "1.008"
$ws.Range("E19").Value = "  +0.69%  "
$ws.Range("D20").Value = "26.951.07"
$ws.Range("E20").Value = "  -1.51%  "
$ws.Range("D21").Value = "14.52"
$ws.Range("E21").Value = "  -2.15%  "
$ws.Range("D22").Value = "5.011"
$ws.Range("E22").Value = "  -2.02%  "
$ws.Range("D23").Value = "2.048.11"
$ws.Range("E23").Value = "  -3.78%  "
$ws.Range("E24").Value = "  -0.82%  "
$ws.Range("D25").Value = "1.985"
$ws.Range("E25").Value = "  +5.28%  "
$ws.Range("D26").Value = "151.35"
$ws.Range("E26").Value = "  -0.49%  "
$ws.Range("D27").Value = "18.20"
$ws.Range("E27").Value = "  -0.95%  "
$ws.Range("D28").Value = "1.970"
$ws.Range("E28").Value = "  -4.93%  "
$ws.Range("D29").Value = "114.02"
$ws.Range("E29").Value = "  -1.65%  "
$ws.Range("D30").Value = "4.938"
$ws.Range("E30").Value = "  -3.71%  "
$ws.Range("D31").Value = "0.08802"
$ws.Range("E31").Value = "  -0.59%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "0.7554"
$ws.Range("E32").Value = "  -1.36%  "
$ws.Range("B33").Value = "HuobiToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D33").Value = "3.020"
$ws.Range("E33").Value = "  -0.03%  "
$ws.Range("D34").Value = "4.487"
$ws.Range("E34").Value = "  -0.08%  "
$ws.Range("D35").Value = "1.133"
$ws.Range("E35").Value = "  -3.36%  "
$ws.Range("D36").Value = "2.573"
$ws.Range("E36").Value = "  -1.92%  "
$ws.Range("D37").Value = "1.090"
$ws.Range("E37").Value = "  +1.08%  "
$ws.Range("E38").Value = "  -1.28%  "
$ws.Range("D39").Value = "0.05142"
$ws.Range("E39").Value = "  -1.27%  "
$ws.Range("D40").Value = "2.890"
$ws.Range("E40").Value = "  -3.33%  "
$ws.Range("D41").Value = "6.932"
$ws.Range("E41").Value = "  -1.69%  "
$ws.Range("D42").Value = "0.4976"
$ws.Range("E42").Value = "  -3.56%  "
$ws.Range("D43").Value = "0.1602"
$ws.Range("E43").Value = "  -2.27%  "
$ws.Range("D44").Value = "8.315"
$ws.Range("E44").Value = "  -0.37%  "
$ws.Range("E45").Value = "  -3.09%  "
$ws.Range("E46").Value = "  +0.78%  "
$ws.Range("D47").Value = "10.16"
$ws.Range("E47").Value = "  -1.56%  "
$ws.Range("D48").Value = "102.27"
$ws.Range("E48").Value = "  -1.09%  "
$ws.Range("D49").Value = "1.614"
$ws.Range("E49").Value = "  -2.32%  "
$ws.Range("D50").Value = "0.06125"
$ws.Range("E50").Value = "  -1.56%  "
$ws.Range("D51").Value = "64.46"
$ws.Range("E51").Value = "  -1.85%  "
